# Updates cryptocurrency price (D) and hourly volume change (E) values
# on the "cryptos" worksheet, matching the latest scrape snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.039.72'
$ws.Range("E2").Value = '  -1.52%  '
$ws.Range("D3").Value = '2.468.99'
$ws.Range("E3").Value = '  -1.37%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").Value = "'518.85"
$ws.Range("E5").Value = '  -2.53%  '
$ws.Range("D6").Value = "'133.56"
$ws.Range("E6").Value = '  -1.41%  '
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = "'0.557"
$ws.Range("E8").Value = '  -1.69%  '
$ws.Range("D9").Value = '2.479.72'
$ws.Range("E9").Value = '  -0.91%  '
$ws.Range("D10").Value = "'0.0982"
$ws.Range("E10").Value = '  -3.12%  '
$ws.Range("D12").Value = "'5.30"
$ws.Range("E12").Value = '  -1.83%  '
$ws.Range("D13").Value = "'0.336"
$ws.Range("E13").Value = '  -2.59%  '
$ws.Range("D14").Value = '2.910.41'
$ws.Range("E14").Value = '  -1.29%  '
$ws.Range("D15").Value = '57.980.46'
$ws.Range("E15").Value = '  -1.52%  '
$ws.Range("D16").Value = "'21.96"
$ws.Range("E16").Value = '  -3.47%  '
$ws.Range("E17").Value = '  -2.25%  '
$ws.Range("D18").Value = '2.481.14'
$ws.Range("E18").Value = '  -0.59%  '
$ws.Range("D19").Value = "'10.60"
$ws.Range("E19").Value = '  -3.88%  '
$ws.Range("D20").Value = "'319.48"
$ws.Range("E20").Value = '  -1.30%  '
$ws.Range("D21").Value = "'4.16"
$ws.Range("E21").Value = '  -1.88%  '
$ws.Range("D23").Value = "'5.73"
$ws.Range("E23").Value = '  -3.10%  '
$ws.Range("D24").Value = "'64.59"
$ws.Range("E24").Value = '  -0.53%  '
$ws.Range("D25").Value = "'0.409"
$ws.Range("E25").Value = '  -2.36%  '
$ws.Range("E26").Value = '  -0.30%  '
$ws.Range("E27").Value = '  -1.32%  '
$ws.Range("D28").Value = "'7.34"
$ws.Range("E28").Value = '  -2.13%  '
$ws.Range("D29").Value = '0.0₃0746'
$ws.Range("E29").Value = '  -2.10%  '
$ws.Range("D30").Value = "'169.30"
$ws.Range("E30").Value = '  +0.23%  '
$ws.Range("E31").Value = '  -1.97%  '
$ws.Range("E32").Value = '  -3.14%  '
$ws.Range("E33").Value = '  +3.98%  '
$ws.Range("D35").Value = "'0.997"
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("D36").Value = "'18.04"
$ws.Range("E36").Value = '  -1.60%  '
$ws.Range("D37").Value = "'1.30"
$ws.Range("E37").Value = '  -3.99%  '
$ws.Range("D38").Value = "'3.99"
$ws.Range("E38").Value = '  -0.99%  '
$ws.Range("D39").Value = "'36.58"
$ws.Range("E39").Value = '  -0.70%  '
$ws.Range("D40").Value = "'1.46"
$ws.Range("E40").Value = '  -3.54%  '
$ws.Range("D41").Value = "'0.796"
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("E42").Value = '  +2.90%  '
$ws.Range("D43").Value = "'274.03"
$ws.Range("E43").Value = '  -2.35%  '
$ws.Range("D44").Value = "'3.43"
$ws.Range("E44").Value = '  -4.13%  '
$ws.Range("D45").Value = "'0.595"
$ws.Range("E45").Value = '  -1.09%  '
$ws.Range("D46").Value = "'123.34"
$ws.Range("E46").Value = '  -4.52%  '
$ws.Range("D47").Value = "'0.0908"
$ws.Range("E47").Value = '  -1.59%  '
$ws.Range("D48").Value = "'0.0488"
$ws.Range("E48").Value = '  -2.20%  '
$ws.Range("D49").Value = "'0.0213"
$ws.Range("E49").Value = '  -2.19%  '
$ws.Range("D50").Value = "'16.98"
$ws.Range("E50").Value = '  -1.43%  '
$ws.Range("D51").Value = '1.733.82'
$ws.Range("E51").Value = '  -0.80%  '
